# "Questions List" workbook update
# - Marks several more questions in column C ("Done [yes or no]") as "yes"
# - Clears the status on the two blank separator rows (42/43) back to empty
# - Leaves the last edited cell selected at B63 (matches where the user
#   scrolled to / finished editing)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Done" cell (column C) should now read "yes"
$doneRows = @(8, 9, 10, 11, 12, 13, 17, 18, 44, 45, 47, 51, 56, 57, 58, 59)
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 3).Value = "yes"
}

# Blank separator rows: clear the stray "<->" status back to empty
$clearRows = @(42, 43)
foreach ($r in $clearRows) {
    $ws.Cells.Item($r, 3).Value = ""
}

# Leave the selection where the user ended up while editing
$ws.Range("B63").Select()
